# Daily attendance processing - 2026-01-18 09:59:59
# Rotate the "Recorded By" list in column G so that the last author listed
# moves to the front of the comma-separated list (for every data row that
# has more than one recorded author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $val = $cell.Value2

    if ($val -ne $null) {
        $s = $val.ToString()
        if ($s.Contains(",")) {
            $parts = $s.Split(",")
            for ($i = 0; $i -lt $parts.Length; $i++) {
                $parts[$i] = $parts[$i].Trim()
            }
            if ($parts.Length -gt 1) {
                $rotated = @($parts[$parts.Length - 1]) + $parts[0..($parts.Length - 2)]
                $cell.Value2 = [string]::Join(", ", $rotated)
            }
        }
    }
}
